$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjusted point distribution: swap points between Q1b (B7) and Q3 (B14)
$ws.Range("B7").Value = 1
$ws.Range("B14").Value = 2

# Update the active cell selection to C22
$ws.Range("C22").Select()
